$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab from "Gamma1F-HW45.xpc" to "Gamma1F"
$ws.Name = "Gamma1F"

# Append a new data row (row 16) mirroring the pattern of row 15
$ws.Range("A16").Value = 14
$ws.Range("A16").Style = "Normal"

$ws.Range("B16").Value = $ws.Range("B15").Value2

$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1

# Copy the formatting of row 15 (bold/centered/bordered A & B cells) into row 16
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
